{"js": "// Assign each section heading of the document to a team member.\n//\n// The document body is a flat list of \"Heading1\" paragraphs (section\n// titles) followed by a \"References\" heading and a body paragraph.\n// For every section heading we append \" - <name>\" (a couple use an\n// en dash \"\u2013\" instead of a hyphen, matching the source edit), we fix\n// the \"GANT Chart\" typo to \"GANTT Chart\", and we insert a brand new\n// \"Idea - Marc\" heading right after \"Introduction\".\n\nconst body = context.document.body;\n\n// The document already carries Word's auto-managed \"_GoBack\" bookmark\n// (marks the last edit position) on the \"Device Functionalities\"\n// heading. The authored edit nudges that bookmark into the middle of\n// \"Introduction\" (between \"Introdu\" and \"ction\") since that is where\n// the author's cursor ended up after editing that heading first.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Look a section heading up by its exact (trimmed) text. Using exact\n// equality -- rather than a substring search -- keeps headings whose\n// text is a prefix/substring of another (e.g. \"Main Components\" vs.\n// \"Capabilities of Main Components\") from being confused.\nfunction getParagraphByText(label) {\n  const match = paragraphs.items.find((p) => p.text.trim() === label);\n  if (!match) {\n    throw new Error(\"paragraph not found: \" + label);\n  }\n  return match;\n}\n\nconst introParagraph = getParagraphByText(\"Introduction\");\n\n// Split \"Introduction\" right after \"Introdu\" and drop the relocated\n// \"_GoBack\" bookmark there, then append the assignee note.\nconst introSplit = introParagraph.search(\"Introdu\", { matchCase: true });\nintroSplit.load(\"items\");\nawait context.sync();\n\nconst splitPoint = introSplit.items[0].getRange(Word.RangeLocation.end);\nsplitPoint.insertBookmark(\"_GoBack\");\nintroParagraph.insertText(\" \\u2013 Marc\", Word.InsertLocation.end);\nawait context.sync();\n\n// Add the brand-new \"Idea - Marc\" heading right after \"Introduction\".\nconst ideaParagraph = introParagraph.insertParagraph(\"Idea - Marc\", Word.InsertLocation.after);\nideaParagraph.styleBuiltIn = Word.BuiltInStyleName.heading1;\nawait context.sync();\n\n// Append the assignee / status note to each remaining heading.\ngetParagraphByText(\"Main Components\").insertText(\" - Marc\", Word.InsertLocation.end);\ngetParagraphByText(\"Capabilities of Main Components\").insertText(\n  \" \\u2013 Will assign once we have clarification of whats required here\",\n  Word.InsertLocation.end\n);\ngetParagraphByText(\"Device Functionalities\").insertText(\" - Robert\", Word.InsertLocation.end);\ngetParagraphByText(\"Hardware Design\").insertText(\" - Patrick\", Word.InsertLocation.end);\ngetParagraphByText(\"Software Design\").insertText(\" - Mark\", Word.InsertLocation.end);\n\n// Keep a handle on the \"GANT Chart\" paragraph before fixing its typo\n// so we can still find it afterwards.\nconst gantParagraph = getParagraphByText(\"GANT Chart\");\nawait context.sync();\n\nconst gantTypo = gantParagraph.search(\"GANT Chart\", { matchCase: true });\ngantTypo.load(\"items\");\nawait context.sync();\ngantTypo.items[0].insertText(\"GANTT Chart\", Word.InsertLocation.replace);\ngantParagraph.insertText(\" - All\", Word.InsertLocation.end);\n\ngetParagraphByText(\"Components Price List\").insertText(\" - Mark\", Word.InsertLocation.end);\ngetParagraphByText(\"Current Progress\").insertText(\n  \" \\u2013 Last weekly report\",\n  Word.InsertLocation.end\n);\ngetParagraphByText(\"Conclusion\").insertText(\" - All\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Assign each section heading of the document to a team member.\n#\n# The document body is a flat list of \"Heading1\" paragraphs (section\n# titles) followed by a \"References\" heading and a body paragraph.\n# For every section heading we append \" - <name>\" (a couple use an en\n# dash \"-\" instead of a hyphen, matching the source edit), we fix the\n# \"GANT Chart\" typo to \"GANTT Chart\", and we insert a brand new\n# \"Idea - Marc\" heading right after \"Introduction\".\n\nfunction Get-ParaByText($doc, $label) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.Trim() -eq $label) {\n            return $p\n        }\n    }\n    throw \"paragraph not found: $label\"\n}\n\n$d = $word.ActiveDocument\n\n# The document already carries Word's auto-managed \"_GoBack\" bookmark\n# (marks the last edit position) at the end of the \"Device\n# Functionalities\" heading. The authored edit nudges that bookmark\n# into the middle of \"Introduction\" (between \"Introdu\" and \"ction\"),\n# since that is where the author's cursor ended up after editing that\n# heading first. Re-adding a bookmark with the same name moves it\n# instead of creating a duplicate.\n$introPara = Get-ParaByText $d \"Introduction\"\n$introStart = $introPara.Range.Start\n$splitRange = $d.Range($introStart + 7, $introStart + 7)\n$d.Bookmarks.Add(\"_GoBack\", $splitRange) | Out-Null\n\n# Insert the brand-new \"Idea - Marc\" heading right after \"Introduction\"\n# (keep a reference to the paragraph object so we can still reach it\n# after its text changes below).\n$introPara.Range.InsertParagraphAfter() | Out-Null\n$ideaPara = $introPara.Next()\n$ideaPara.Range.Text = \"Idea - Marc\"\n\n# Append the assignee note to \"Introduction\" (do this last, since its\n# text no longer matches \"Introduction\" exactly afterwards).\n$introPara.Range.InsertAfter(\" \" + [char]0x2013 + \" Marc\") | Out-Null\n\n# Append the assignee / status note to each remaining heading.\n(Get-ParaByText $d \"Main Components\").Range.InsertAfter(\" - Marc\") | Out-Null\n(Get-ParaByText $d \"Capabilities of Main Components\").Range.InsertAfter(\n    \" \" + [char]0x2013 + \" Will assign once we have clarification of whats required here\") | Out-Null\n(Get-ParaByText $d \"Device Functionalities\").Range.InsertAfter(\" - Robert\") | Out-Null\n(Get-ParaByText $d \"Hardware Design\").Range.InsertAfter(\" - Patrick\") | Out-Null\n(Get-ParaByText $d \"Software Design\").Range.InsertAfter(\" - Mark\") | Out-Null\n\n# Fix the \"GANT Chart\" typo, then append its assignee note.\n$gantPara = Get-ParaByText $d \"GANT Chart\"\n$find = $gantPara.Range.Find\n$find.Execute(\"GANT Chart\", $false, $false, $false, $false, $false, $true, 1, $false, \"GANTT Chart\", 2) | Out-Null\n$gantPara = Get-ParaByText $d \"GANTT Chart\"\n$gantPara.Range.InsertAfter(\" - All\") | Out-Null\n\n(Get-ParaByText $d \"Components Price List\").Range.InsertAfter(\" - Mark\") | Out-Null\n(Get-ParaByText $d \"Current Progress\").Range.InsertAfter(\" \" + [char]0x2013 + \" Last weekly report\") | Out-Null\n(Get-ParaByText $d \"Conclusion\").Range.InsertAfter(\" - All\") | Out-Null\n"}
